# Fix the issue when getting a resource after a post:
#  - The Billing sheet's sample "id" value was stale (155); update it to 11.
#  - Re-point the active tab / cell selection so the workbook opens on the
#    Billing sheet at A2 instead of ProductDesc at C2.

$wb = $excel.ActiveWorkbook

$wsBilling     = $wb.Worksheets.Item("Billing")
$wsProductDesc = $wb.Worksheets.Item("ProductDesc")

# Move ProductDesc's own selection back to A2 (it was C2) *before* switching
# the active sheet away from it, so this only updates its stored selection
# without making it the active tab.
$wsProductDesc.Range("A2").Select()

# Correct the stale sample id on the Billing sheet.
$wsBilling.Range("A2").Value = 11

# Make Billing the active sheet with A2 selected (previously ProductDesc was
# the active/selected tab, and Billing's own stored selection was A24).
$wsBilling.Activate()
$wsBilling.Range("A2").Select()
